$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1102
$ws.Range("J31").Value = 7500
$ws.Range("L31").Value = 22500
$ws.Range("N31").Value = -22960

$ws.Range("H64").Value = 3888.0576
$ws.Range("I64").Value = 3773.913
$ws.Range("J64").Value = 3978.5862
$ws.Range("K64").Value = 3773.913
$ws.Range("L64").Value = 3978.5862
$ws.Range("M64").Value = -3525.913
$ws.Range("N64").Value = -4474.5862

$ws.Range("H67").Value = 3888.0576
$ws.Range("I67").Value = 3773.913
$ws.Range("J67").Value = 3978.5862
$ws.Range("K67").Value = 3773.913
$ws.Range("L67").Value = 3978.5862
$ws.Range("M67").Value = -2915.913
$ws.Range("N67").Value = -5694.5862

$ws.Range("H76").Value = 2868.6875
$ws.Range("I76").Value = 2188.875
$ws.Range("J76").Value = 3548.5
$ws.Range("K76").Value = 2188.875
$ws.Range("L76").Value = 3548.5
$ws.Range("M76").Value = -1873.875
$ws.Range("N76").Value = -4178.5

$ws.Range("H79").Value = 2868.6875
$ws.Range("I79").Value = 2188.875
$ws.Range("J79").Value = 3548.5
$ws.Range("K79").Value = 2188.875
$ws.Range("L79").Value = 3548.5
$ws.Range("M79").Value = -1096.875
$ws.Range("N79").Value = -5732.5

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 25001244
$ws.Range("I26").Value = 25001244
$ws.Range("K26").Value = 25001244
$ws.Range("M26").Value = -25000914

$ws.Range("H39").Value = 4999
$ws.Range("I39").Value = 4999
$ws.Range("K39").Value = 4999
$ws.Range("M39").Value = -4479

$ws.Range("H61").Value = 2238.5293
$ws.Range("I61").Value = 1783.4814
$ws.Range("K61").Value = 1783.4814
$ws.Range("M61").Value = -1571.4814

$ws.Range("H102").Value = 2114.4443
$ws.Range("I102").Value = 1684
$ws.Range("J102").Value = 2652.5
$ws.Range("K102").Value = 1684
$ws.Range("L102").Value = 2652.5
$ws.Range("M102").Value = -62
$ws.Range("N102").Value = -5896.5

$ws.Range("H110").Value = 2210.1428
$ws.Range("I110").Value = 1526.6666
$ws.Range("J110").Value = 3918.8333
$ws.Range("K110").Value = 1526.6666
$ws.Range("L110").Value = 3918.8333
$ws.Range("M110").Value = 518.3334
$ws.Range("N110").Value = -8008.8333

$ws.Range("H136").Value = 2238.5293
$ws.Range("I136").Value = 1783.4814
$ws.Range("K136").Value = 5350.4442
$ws.Range("M136").Value = -2800.4442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1514.3334
$ws.Range("I36").Value = 1514.3334
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1514.3334
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -980.3334
$ws.Range("N36").Value = ""

$ws.Range("H57").Value = 48000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 48000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 48000
$ws.Range("M57").Value = ""
$ws.Range("N57").Value = -49440

$ws.Range("H99").Value = 2010
$ws.Range("I99").Value = 1470.9
$ws.Range("K99").Value = 1470.9
$ws.Range("M99").Value = 27.09999999999991

$ws.Range("H105").Value = 2959.7307
$ws.Range("I105").Value = 2502.5
$ws.Range("K105").Value = 2502.5
$ws.Range("M105").Value = -755.5

$ws.Range("H107").Value = 3514.6365
$ws.Range("I107").Value = 3514.6365
$ws.Range("K107").Value = 3514.6365
$ws.Range("M107").Value = -1594.6365

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

$ws.Range("H123").Value = 25780
$ws.Range("J123").Value = 25780
$ws.Range("L123").Value = 25780
$ws.Range("N123").Value = -35580

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""

$ws.Range("H125").Value = 40000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 40000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 40000
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -49840

$ws.Range("H127").Value = 19800
$ws.Range("J127").Value = 19800
$ws.Range("L127").Value = 19800
$ws.Range("N127").Value = -29720

$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

$ws.Range("H130").Value = 42780
$ws.Range("J130").Value = 42780
$ws.Range("L130").Value = 42780
$ws.Range("N130").Value = -52820

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""

$ws.Range("H132").Value = 42000
$ws.Range("J132").Value = 42000
$ws.Range("L132").Value = 42000
$ws.Range("N132").Value = -52120

$ws.Range("H133").Value = 50390
$ws.Range("J133").Value = 50390
$ws.Range("L133").Value = 50390
$ws.Range("N133").Value = -60510

$ws.Range("H134").Value = 1629.6562
$ws.Range("I134").Value = 1244.1072
$ws.Range("J134").Value = 4328.5
$ws.Range("K134").Value = 3732.3216
$ws.Range("L134").Value = 12985.5
$ws.Range("M134").Value = -1197.3216
$ws.Range("N134").Value = -18055.5

$ws.Range("H135").Value = 25000
$ws.Range("I135").Value = 30000
$ws.Range("K135").Value = 30000
$ws.Range("M135").Value = -24930

$ws.Range("H136").Value = 48000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 48000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 48000
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -58200

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws.Range("H138").Value = 82999.39999999999
$ws.Range("J138").Value = 82999.39999999999
$ws.Range("L138").Value = 82999.39999999999
$ws.Range("N138").Value = -93279.39999999999

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws.Range("H141").Value = 34780
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 34780
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 34780
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -45140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2996.309
$ws.Range("I31").Value = 1663.7931
$ws.Range("J31").Value = 4482.577
$ws.Range("K31").Value = 1663.7931
$ws.Range("L31").Value = 4482.577
$ws.Range("M31").Value = -1368.7931
$ws.Range("N31").Value = -5072.577

$ws.Range("H34").Value = 2996.309
$ws.Range("I34").Value = 1663.7931
$ws.Range("J34").Value = 4482.577
$ws.Range("K34").Value = 1663.7931
$ws.Range("L34").Value = 4482.577
$ws.Range("M34").Value = -1461.7931
$ws.Range("N34").Value = -4886.577

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1996.6666
$ws.Range("J41").Value = 1996.6666
$ws.Range("L41").Value = 5989.9998
$ws.Range("N41").Value = -6665.9998

$ws.Range("H131").Value = 753
$ws.Range("I131").Value = 417.125
$ws.Range("J131").Value = 894.4211
$ws.Range("K131").Value = 1251.375
$ws.Range("L131").Value = 2683.2633
$ws.Range("M131").Value = 3788.625
$ws.Range("N131").Value = -12763.2633

$ws.Range("H132").Value = 1982.9
$ws.Range("I132").Value = 502
$ws.Range("K132").Value = 4518
$ws.Range("M132").Value = -1988
